$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by
# Range.Value (e.g. "246.77", "42.30", "4.50") are written with a leading
# apostrophe to force text, then restyled to 'Normal' so no stray number
# format / quote-prefix style is left behind on the cell.

$ws.Range('D2').Value = '42.214.32'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').Value = '2.244.09'
$ws.Range('E3').Value = '  -2.35%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'246.77"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.39%  '
$ws.Range('D6').Value = "'0.623"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.49%  '
$ws.Range('D7').Value = "'76.73"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.88%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = "'0.613"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.34%  '
$ws.Range('D10').Value = "'42.30"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.33%  '
$ws.Range('D11').Value = "'0.0949"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.80%  '
$ws.Range('D12').Value = "'7.07"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -7.14%  '
$ws.Range('D13').Value = "'0.102"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.29%  '
$ws.Range('D14').Value = '2.580.92'
$ws.Range('E14').Value = '  -2.34%  '
$ws.Range('D15').Value = "'14.71"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.78%  '
$ws.Range('D16').Value = "'0.855"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.86%  '
$ws.Range('D17').Value = '2.236.85'
$ws.Range('E17').Value = '  -2.65%  '
$ws.Range('D18').Value = '42.046.69'
$ws.Range('E18').Value = '  -2.00%  '
$ws.Range('D19').Value = '0.0₃0979'
$ws.Range('E19').Value = '  -3.71%  '
$ws.Range('D20').Value = "'71.83"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.40%  '
$ws.Range('D21').Value = "'6.08"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.76%  '
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('D23').Value = "'229.92"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.68%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').Value = "'11.30"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.54%  '
$ws.Range('E26').Value = '  -7.52%  '
$ws.Range('D27').Value = "'2.29"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.32%  '
$ws.Range('D28').Value = "'7.37"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +16.40%  '
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('D30').Value = "'169.68"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.13%  '
$ws.Range('E31').Value = '  -2.77%  '
$ws.Range('D32').Value = "'0.0834"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.26%  '
$ws.Range('D33').Value = "'32.32"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.86%  '
$ws.Range('E34').Value = '  -5.46%  '
$ws.Range('D35').Value = "'0.124"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.76%  '
$ws.Range('D36').Value = "'4.50"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.85%  '
$ws.Range('D37').Value = "'4.95"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.13%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = "'0.0301"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.44%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = "'14.16"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.95%  '
$ws.Range('B40').Value = 'THORChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D40').Value = "'5.88"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').Value = "'2.18"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.01%  '
$ws.Range('D42').Value = "'112.56"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.80%  '
$ws.Range('D43').Value = "'0.202"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.83%  '
$ws.Range('D44').Value = "'60.72"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.67%  '
$ws.Range('D45').Value = "'8.66"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.19%  '
$ws.Range('D46').Value = "'0.0988"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.99%  '
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').Value = "'1.12"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.26%  '
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('D50').Value = "'4.27"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -12.40%  '
$ws.Range('D51').Value = "'0.435"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +13.22%  '
